# Merged with PV sequence manager code
# Fix variable generator script
#
# Adds new global variables to the "Constants" sheet (PV state-machine /
# PV-change-sequence bookkeeping variables) that were introduced by the
# PV sequence manager merge.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Step 1: variable name (B), array size (C), type (D) for each new row ---
$ws.Range("B29").Value = "main_current_state"
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = "ARRAY [4] OF INT"

$ws.Range("B30").Value = "main_next_state"
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = "ARRAY [4] OF INT"

$ws.Range("B31").Value = "transition_current_state"
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = "ARRAY [4] OF INT"

$ws.Range("B32").Value = "transition_next_state"
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = "ARRAY [4] OF INT"

$ws.Range("B33").Value = "PVChangeMode"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = "INT"

$ws.Range("B34").Value = "delay_counter"
$ws.Range("C34").Value = 4
$ws.Range("D34").Value = "ARRAY [4] OF WORD"

$ws.Range("B35").Value = "lastPVvalue"
$ws.Range("C35").Value = 20
$ws.Range("D35").Value = "ARRAY [10] OF REAL"

# --- Step 2: default value (E) for each new row ---
$ws.Range("E29").Value = "[4(0)]"
$ws.Range("E30").Value = "[4(0)]"
$ws.Range("E31").Value = "[4(0)]"
$ws.Range("E32").Value = "[4(0)]"
$ws.Range("E33").Value = 0
$ws.Range("E34").Value = "[4(0)]"

# --- Step 3: description (G) for each new row (row 33 has none) ---
$ws.Range("G29").Value = "Main state machine current state (for each instance)"
$ws.Range("G30").Value = "Main state machine next state (for each instance)"
$ws.Range("G31").Value = "PV transition state machine current state (for each instance)"
$ws.Range("G32").Value = "PV transition state machine next state (for each instance)"
$ws.Range("G34").Value = "Counter for delays"
$ws.Range("G35").Value = "Store last PV value for comparison whether to run PV change sequence or not"

# --- Step 4: final default value for the REAL array, added last ---
$ws.Range("E35").Value = "[10(0.0)]"

# Update the view: scroll down toward the newly added rows, select the
# next empty row under the new data (mirrors author's saved cursor state).
[void]$ws.Range("B36").Select()
